# Updated with BRIEF - DELGADILLO BRIEF.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert a new row at position 18 ("InputFolderPath" setting), pushing the
# existing rows (the blank spacer + all settings below it) down by one.
$ws.Rows("18:18").Insert()

$ws.Range("A18").Value = "InputFolderPath"
$ws.Range("B18").Value = "C:\Users\ewong-c\Documents\UiPath\3DCA_Workflow\Data\Input\"

# Keep the inserted row (and the blank spacer row that got pushed to 19)
# consistent with the sheet's normal row height.
$ws.Rows("18:18").RowHeight = 14.25
$ws.Rows("19:19").RowHeight = 14.25

# The "Delgadillo_Headers" value (now on row 50 after the insert) gets four
# additional pattern entries appended.
$ws.Range("B50").Value = "Certificate of Compliance,Word Count Certification,Word Count Certificate,Certificate of Length,Certificate of Counsel,Certification of Word Count,Brief Format Certification,Certificate of Word Count,Declaration Regarding Word Count,Statement of Word Count,Certification of Compliance,Word Count@Statement of Appealability,Statement of Appellate Jurisdiction,Statement of Jurisdiction,Statement of Appelability,Appealability@Table of Contents,Topical Index@Table of Authorities,Table of Authority,Table of Cases and Authorities,Table of Cases@Statement of Facts,Statement of the Facts,Statement of the Case,Summary of Alleged Facts,Statements of Case and Procedural Facts@Proof of Service,Certificate of Service,Proof of Service by Mail and Email,Proof of Service by Mail,Proof of Service by Mail/Electronic Service,Declaration of Service,Declaration of Service by Mail,Proof of Service by Mail and Electronically,Brief Format Certification Pursuant to California Rules of Court 8.360,Declaration of Electronic Service and Service by U.S. Mail,Attorney's Certificate of Electronic Service and Service by Mail,Attorney's Certificate of Service,@Argument,Arguments@Declaration of Appellate Counsel,Declaration of Counsel,DECLARATION OF,COUNSEL DECLARATION,Certification of Appellate Counsel,Certificate of Appellate Counsel,Attest of Appellate Counsel,Attestation of Appellate Counsel@court may treat the appeal as abandoned and dismiss the appeal,court may dismiss the appeal without issuing an opinion,may result in the court dismissing the appeal as abandoned,the appeal may be dismissed as abandoned,court very likely would dismiss the appeal as abandoned,no brief or letter is filed the appeal may be dismissed,court may treat the appeal as abandoned and dismiss,the court may dismiss the matter,the Court may dismiss the appeal as abandoned,if the appeal is dismissed as abandoned,Court may treat his appeal as abandoned and dismiss,Court may treat her appeal as abandoned and dismiss,Court may treat their appeal as abandoned and dismiss,Court may treat they appeal as abandoned and dismiss"
